$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 468
$ws1.Range("F4").Value = 7843
$ws1.Range("F6").Value = 212
$ws1.Range("F10").Value = 456
$ws1.Range("F13").Value = 443
$ws1.Range("F14").Value = 65
$ws1.Range("F15").Value = 68
$ws1.Range("F17").Value = 5725
$ws1.Range("G17").Value = 44.1
$ws1.Range("F19").Value = 240
$ws1.Range("F20").Value = 1407

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 468
$ws4.Range("F4").Value = 7843
$ws4.Range("F6").Value = 212
$ws4.Range("F10").Value = 456
$ws4.Range("F13").Value = 443
$ws4.Range("F14").Value = 65
$ws4.Range("F15").Value = 68
$ws4.Range("F18").Value = 5725
$ws4.Range("G18").Value = 44.1
$ws4.Range("F21").Value = 240
$ws4.Range("F22").Value = 1407
